# Update TPM-derived NATMI ligand/receptor statistics (Efnb2-Ephb3)
# with refreshed values from the rerun of the scripts ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.94436433333333
$ws.Range("H2").Value = 110.833093
$ws.Range("I2").Value = 0.8328964975864823
$ws.Range("J2").Value = 0.8328964975864824
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01631833333333333
$ws.Range("N2").Value = 0.048955
$ws.Range("O2").Value = 0.001076315602073535
$ws.Range("P2").Value = 0.001076315602073535
$ws.Range("Q2").Value = 0.6028704519794444
$ws.Range("R2").Value = 5.425834067815
$ws.Range("S2").Value = 0.0008964594952647334
$ws.Range("T2").Value = 0.0008964594952647334
$ws.Range("G3").Value = 36.94436433333333
$ws.Range("H3").Value = 110.833093
$ws.Range("I3").Value = 0.8328964975864823
$ws.Range("J3").Value = 0.8328964975864824
$ws.Range("O3").Value = 0.7730166590262294
$ws.Range("P3").Value = 0.7730166590262293
$ws.Range("Q3").Value = 432.9853638811633
$ws.Range("R3").Value = 3896.86827493047
$ws.Range("S3").Value = 0.6438428678789505
$ws.Range("T3").Value = 0.6438428678789505
$ws.Range("G4").Value = 36.94436433333333
$ws.Range("H4").Value = 110.833093
$ws.Range("I4").Value = 0.8328964975864823
$ws.Range("J4").Value = 0.8328964975864824
$ws.Range("M4").Value = 3.425042
$ws.Range("N4").Value = 10.275126
$ws.Range("O4").Value = 0.2259070253716972
$ws.Range("P4").Value = 0.2259070253716972
$ws.Range("Q4").Value = 126.5359995049687
$ws.Range("R4").Value = 1138.823995544718
$ws.Range("S4").Value = 0.1881571702122672
$ws.Range("T4").Value = 0.1881571702122672
$ws.Range("I5").Value = 0.07608399754092349
$ws.Range("J5").Value = 0.07608399754092349
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01631833333333333
$ws.Range("N5").Value = 0.048955
$ws.Range("O5").Value = 0.001076315602073535
$ws.Range("P5").Value = 0.001076315602073535
$ws.Range("Q5").Value = 0.05507142138166667
$ws.Range("R5").Value = 0.495642792435
$ws.Range("S5").Value = 0.00008189039362142043
$ws.Range("T5").Value = 0.00008189039362142041
$ws.Range("I6").Value = 0.07608399754092349
$ws.Range("J6").Value = 0.07608399754092349
$ws.Range("O6").Value = 0.7730166590262294
$ws.Range("P6").Value = 0.7730166590262293
$ws.Range("S6").Value = 0.05881419758444453
$ws.Range("T6").Value = 0.05881419758444452
$ws.Range("I7").Value = 0.07608399754092349
$ws.Range("J7").Value = 0.07608399754092349
$ws.Range("M7").Value = 3.425042
$ws.Range("N7").Value = 10.275126
$ws.Range("O7").Value = 0.2259070253716972
$ws.Range("P7").Value = 0.2259070253716972
$ws.Range("Q7").Value = 11.558896817398
$ws.Range("R7").Value = 104.030071356582
$ws.Range("S7").Value = 0.01718790956285755
$ws.Range("T7").Value = 0.01718790956285755
$ws.Range("G8").Value = 4.037305666666668
$ws.Range("H8").Value = 12.111917
$ws.Range("I8").Value = 0.09101950487259411
$ws.Range("J8").Value = 0.09101950487259411
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01631833333333333
$ws.Range("N8").Value = 0.048955
$ws.Range("O8").Value = 0.001076315602073535
$ws.Range("P8").Value = 0.001076315602073535
$ws.Range("Q8").Value = 0.06588209963722223
$ws.Range("R8").Value = 0.5929388967350001
$ws.Range("S8").Value = 0.0000979657131873812
$ws.Range("T8").Value = 0.00009796571318738117
$ws.Range("G9").Value = 4.037305666666668
$ws.Range("H9").Value = 12.111917
$ws.Range("I9").Value = 0.09101950487259411
$ws.Range("J9").Value = 0.09101950487259411
$ws.Range("O9").Value = 0.7730166590262294
$ws.Range("P9").Value = 0.7730166590262293
$ws.Range("Q9").Value = 47.31693980193668
$ws.Range("R9").Value = 425.8524582174301
$ws.Range("S9").Value = 0.0703595935628343
$ws.Range("T9").Value = 0.07035959356283429
$ws.Range("G10").Value = 4.037305666666668
$ws.Range("H10").Value = 12.111917
$ws.Range("I10").Value = 0.09101950487259411
$ws.Range("J10").Value = 0.09101950487259411
$ws.Range("M10").Value = 3.425042
$ws.Range("N10").Value = 10.275126
$ws.Range("O10").Value = 0.2259070253716972
$ws.Range("P10").Value = 0.2259070253716972
$ws.Range("Q10").Value = 13.82794147517134
$ws.Range("R10").Value = 124.451473276542
$ws.Range("S10").Value = 0.02056194559657243
$ws.Range("T10").Value = 0.02056194559657243
